$wb = $excel.ActiveWorkbook

# Turn off alerts so sheet deletion doesn't prompt a confirmation dialog
$excel.DisplayAlerts = $false

# Remove the sheets that are no longer part of the template
$wb.Worksheets.Item("MissingParam").Delete()
$wb.Worksheets.Item("Aciclovir").Delete()

$excel.DisplayAlerts = $true

# Rename the remaining sheet from "Global" to "Template"
$ws = $wb.Worksheets.Item("Global")
$ws.Name = "Template"

# Update the selection to reflect where the cursor ended up when the
# workbook was last saved
$ws.Activate()
$ws.Range("A32").Select()

$wb.Save()
